$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '58.358.85'
Set-TextValue "E2" '  -0.91%  '
Set-TextValue "D3" '2.486.04'
Set-TextValue "E3" '  -0.42%  '
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  +0.08%  '
Set-TextValue "D5" '521.41'
Set-TextValue "E5" '  -1.93%  '
Set-TextValue "D6" '135.31'
Set-TextValue "E6" '  +0.25%  '
Set-TextValue "D7" '0.998'
Set-TextValue "E7" '  -0.15%  '
Set-TextValue "E8" '  -1.24%  '
Set-TextValue "D9" '2.504.52'
Set-TextValue "E9" '  +0.26%  '
Set-TextValue "D10" '0.0992'
Set-TextValue "E10" '  -2.37%  '
Set-TextValue "D11" '0.157'
Set-TextValue "E11" '  -0.70%  '
Set-TextValue "E12" '  -0.61%  '
Set-TextValue "E13" '  -1.67%  '
Set-TextValue "D14" '2.928.74'
Set-TextValue "E14" '  -0.32%  '
Set-TextValue "D15" '58.286.73'
Set-TextValue "E15" '  -0.88%  '
Set-TextValue "E16" '  -1.84%  '
Set-TextValue "E17" '  -1.66%  '
Set-TextValue "D18" '2.497.86'
Set-TextValue "E18" '  -0.21%  '
Set-TextValue "D19" '10.71'
Set-TextValue "E19" '  -2.78%  '
Set-TextValue "E20" '  -0.98%  '
Set-TextValue "D21" '322.02'
Set-TextValue "E21" '  -0.25%  '
Set-TextValue "D23" '5.77'
Set-TextValue "E23" '  -2.80%  '
Set-TextValue "D24" '64.53'
Set-TextValue "E24" '  -0.68%  '
Set-TextValue "E25" '  -1.66%  '
Set-TextValue "B26" 'Binance-PegBSC-USD'
Set-TextValue "C26" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D26" '0.997'
Set-TextValue "E26" '  -0.28%  '
Set-TextValue "B27" 'Kaspa'
Set-TextValue "C27" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D27" '0.162'
Set-TextValue "E27" '  -1.20%  '
Set-TextValue "E28" '  -1.03%  '
Set-TextValue "D29" '0.0₃0752'
Set-TextValue "D30" '169.94'
Set-TextValue "E30" '  -0.07%  '
Set-TextValue "B31" 'PancakeSwap'
Set-TextValue "C31" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D31" '1.70'
Set-TextValue "E31" '  -2.19%  '
Set-TextValue "B32" 'Aptos'
Set-TextValue "C32" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D32" '6.34'
Set-TextValue "E32" '  -1.32%  '
Set-TextValue "E33" '  +2.53%  '
Set-TextValue "E34" '  -0.06%  '
Set-TextValue "D35" '0.996'
Set-TextValue "E35" '  -0.27%  '
Set-TextValue "E36" '  -0.80%  '
Set-TextValue "E37" '  -0.49%  '
Set-TextValue "D38" '4.05'
Set-TextValue "E38" '  +0.21%  '
Set-TextValue "D39" '36.70'
Set-TextValue "E39" '  -0.12%  '
Set-TextValue "E40" '  -2.71%  '
Set-TextValue "D41" '0.802'
Set-TextValue "E41" '  +0.34%  '
Set-TextValue "D42" '5.22'
Set-TextValue "E42" '  +4.45%  '
Set-TextValue "D43" '279.24'
Set-TextValue "E43" '  -0.50%  '
Set-TextValue "E44" '  -2.49%  '
Set-TextValue "D45" '0.601'
Set-TextValue "E45" '  +0.42%  '
Set-TextValue "D46" '124.37'
Set-TextValue "E46" '  -4.02%  '
Set-TextValue "D47" '0.0911'
Set-TextValue "E47" '  -1.39%  '
Set-TextValue "E48" '  -0.95%  '
Set-TextValue "E49" '  -1.34%  '
Set-TextValue "D50" '17.16'
Set-TextValue "E50" '  -0.27%  '
Set-TextValue "D51" '1.743.83'
Set-TextValue "E51" '  -0.36%  '
